$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab
$ws.Name = "Through 2022-06-06"

# Update header text for first month column
$ws.Range("B1").Value = "June 2022 (through June 06)"

# Simple single-cell value updates
$ws.Range("B2").Value = 3
$ws.Range("H2").Value = 5
$ws.Range("N2").Value = 5
$ws.Range("N3").Value = 3
$ws.Range("H6").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("H9").Value = 1
$ws.Range("N10").Value = 4
$ws.Range("N29").Value = 1
$ws.Range("H48").Value = 1

# Prime row 97 (new row) with the style used by column-A header cells, copied from row 96
$ws.Range("A96").Copy()
$ws.Range("A97").PasteSpecial(-4122)

# Shift rows 61-96 down to 62-97 (process bottom-up), clearing destination first
$ws.Range("B97:AW97").ClearContents()
$ws.Range("A97").Value = 'Wrigleyville'
$ws.Range("E97").Value = 1
$ws.Range("N97").Value = 1
$ws.Range("AB97").Value = 1
$ws.Range("AK97").Value = 1
$ws.Range("B96:AW96").ClearContents()
$ws.Range("A96").Value = 'Wicker Park'
$ws.Range("D96").Value = 1
$ws.Range("E96").Value = 3
$ws.Range("F96").Value = 3
$ws.Range("G96").Value = 1
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 2
$ws.Range("L96").Value = 2
$ws.Range("M96").Value = 7
$ws.Range("O96").Value = 1
$ws.Range("P96").Value = 2
$ws.Range("Q96").Value = 1
$ws.Range("R96").Value = 1
$ws.Range("S96").Value = 2
$ws.Range("U96").Value = 2
$ws.Range("V96").Value = 1
$ws.Range("X96").Value = 1
$ws.Range("AA96").Value = 1
$ws.Range("AF96").Value = 1
$ws.Range("AH96").Value = 1
$ws.Range("AJ96").Value = 1
$ws.Range("AN96").Value = 1
$ws.Range("AT96").Value = 1
$ws.Range("B95:AW95").ClearContents()
$ws.Range("A95").Value = 'West Town'
$ws.Range("D95").Value = 1
$ws.Range("E95").Value = 2
$ws.Range("F95").Value = 4
$ws.Range("G95").Value = 4
$ws.Range("I95").Value = 3
$ws.Range("K95").Value = 3
$ws.Range("L95").Value = 6
$ws.Range("M95").Value = 1
$ws.Range("O95").Value = 1
$ws.Range("P95").Value = 1
$ws.Range("Q95").Value = 1
$ws.Range("R95").Value = 1
$ws.Range("S95").Value = 1
$ws.Range("U95").Value = 1
$ws.Range("V95").Value = 1
$ws.Range("W95").Value = 1
$ws.Range("X95").Value = 1
$ws.Range("AA95").Value = 3
$ws.Range("AB95").Value = 2
$ws.Range("AG95").Value = 1
$ws.Range("AH95").Value = 2
$ws.Range("AI95").Value = 1
$ws.Range("AJ95").Value = 1
$ws.Range("AK95").Value = 2
$ws.Range("AM95").Value = 1
$ws.Range("B94:AW94").ClearContents()
$ws.Range("A94").Value = 'West Pullman'
$ws.Range("B94").Value = 1
$ws.Range("D94").Value = 1
$ws.Range("F94").Value = 3
$ws.Range("G94").Value = 1
$ws.Range("I94").Value = 2
$ws.Range("J94").Value = 2
$ws.Range("L94").Value = 1
$ws.Range("N94").Value = 1
$ws.Range("O94").Value = 3
$ws.Range("P94").Value = 2
$ws.Range("R94").Value = 1
$ws.Range("S94").Value = 2
$ws.Range("T94").Value = 1
$ws.Range("V94").Value = 2
$ws.Range("AC94").Value = 1
$ws.Range("AD94").Value = 1
$ws.Range("AE94").Value = 1
$ws.Range("AM94").Value = 2
$ws.Range("AN94").Value = 3
$ws.Range("AO94").Value = 3
$ws.Range("B93:AW93").ClearContents()
$ws.Range("A93").Value = 'West Elsdon'
$ws.Range("E93").Value = 1
$ws.Range("H93").Value = 1
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = 1
$ws.Range("L93").Value = 1
$ws.Range("AD93").Value = 1
$ws.Range("AK93").Value = 1
$ws.Range("AM93").Value = 1
$ws.Range("AU93").Value = 1
$ws.Range("B92:AW92").ClearContents()
$ws.Range("A92").Value = 'Washington Heights'
$ws.Range("B92").Value = 1
$ws.Range("D92").Value = 3
$ws.Range("E92").Value = 3
$ws.Range("F92").Value = 3
$ws.Range("G92").Value = 2
$ws.Range("H92").Value = 1
$ws.Range("J92").Value = 1
$ws.Range("L92").Value = 1
$ws.Range("M92").Value = 4
$ws.Range("N92").Value = 1
$ws.Range("O92").Value = 1
$ws.Range("P92").Value = 2
$ws.Range("Q92").Value = 3
$ws.Range("S92").Value = 1
$ws.Range("X92").Value = 1
$ws.Range("Y92").Value = 2
$ws.Range("AA92").Value = 2
$ws.Range("AC92").Value = 2
$ws.Range("AD92").Value = 1
$ws.Range("AG92").Value = 3
$ws.Range("AH92").Value = 1
$ws.Range("AI92").Value = 1
$ws.Range("AJ92").Value = 1
$ws.Range("AK92").Value = 2
$ws.Range("AM92").Value = 2
$ws.Range("AN92").Value = 1
$ws.Range("AP92").Value = 2
$ws.Range("AQ92").Value = 2
$ws.Range("AS92").Value = 1
$ws.Range("B91:AW91").ClearContents()
$ws.Range("A91").Value = 'Uptown'
$ws.Range("D91").Value = 1
$ws.Range("E91").Value = 2
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 1
$ws.Range("I91").Value = 2
$ws.Range("J91").Value = 3
$ws.Range("L91").Value = 2
$ws.Range("M91").Value = 1
$ws.Range("O91").Value = 1
$ws.Range("U91").Value = 1
$ws.Range("V91").Value = 1
$ws.Range("AB91").Value = 1
$ws.Range("AF91").Value = 1
$ws.Range("AG91").Value = 1
$ws.Range("AI91").Value = 1
$ws.Range("AM91").Value = 1
$ws.Range("AQ91").Value = 1
$ws.Range("B90:AW90").ClearContents()
$ws.Range("A90").Value = 'Streeterville'
$ws.Range("H90").Value = 1
$ws.Range("I90").Value = 1
$ws.Range("J90").Value = 1
$ws.Range("O90").Value = 1
$ws.Range("Q90").Value = 1
$ws.Range("AD90").Value = 1
$ws.Range("AH90").Value = 1
$ws.Range("AI90").Value = 1
$ws.Range("B89:AW89").ClearContents()
$ws.Range("A89").Value = 'South Deering'
$ws.Range("D89").Value = 1
$ws.Range("E89").Value = 1
$ws.Range("G89").Value = 2
$ws.Range("I89").Value = 1
$ws.Range("J89").Value = 2
$ws.Range("K89").Value = 3
$ws.Range("M89").Value = 1
$ws.Range("N89").Value = 1
$ws.Range("Q89").Value = 1
$ws.Range("U89").Value = 1
$ws.Range("Y89").Value = 1
$ws.Range("AA89").Value = 1
$ws.Range("AG89").Value = 1
$ws.Range("AN89").Value = 1
$ws.Range("AP89").Value = 1
$ws.Range("AS89").Value = 1
$ws.Range("AT89").Value = 1
$ws.Range("B88:AW88").ClearContents()
$ws.Range("A88").Value = 'Sheffield & DePaul'
$ws.Range("E88").Value = 1
$ws.Range("M88").Value = 1
$ws.Range("R88").Value = 1
$ws.Range("AB88").Value = 1
$ws.Range("B87:AW87").ClearContents()
$ws.Range("A87").Value = 'Sauganash,Forest Glen'
$ws.Range("M87").Value = 2
$ws.Range("B86:AW86").ClearContents()
$ws.Range("A86").Value = 'Rush & Division'
$ws.Range("R86").Value = 1
$ws.Range("AD86").Value = 1
$ws.Range("AG86").Value = 1
$ws.Range("B85:AW85").ClearContents()
$ws.Range("A85").Value = 'River North'
$ws.Range("D85").Value = 1
$ws.Range("F85").Value = 2
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 1
$ws.Range("I85").Value = 1
$ws.Range("K85").Value = 1
$ws.Range("M85").Value = 3
$ws.Range("O85").Value = 2
$ws.Range("P85").Value = 1
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = 1
$ws.Range("S85").Value = 1
$ws.Range("W85").Value = 1
$ws.Range("Y85").Value = 1
$ws.Range("AB85").Value = 1
$ws.Range("AC85").Value = 1
$ws.Range("AG85").Value = 1
$ws.Range("AH85").Value = 1
$ws.Range("AI85").Value = 2
$ws.Range("AJ85").Value = 2
$ws.Range("AP85").Value = 1
$ws.Range("AT85").Value = 1
$ws.Range("B84:AW84").ClearContents()
$ws.Range("A84").Value = 'Pullman'
$ws.Range("D84").Value = 2
$ws.Range("J84").Value = 2
$ws.Range("O84").Value = 1
$ws.Range("AV84").Value = 1
$ws.Range("B83:AW83").ClearContents()
$ws.Range("A83").Value = 'Printers Row'
$ws.Range("J83").Value = 2
$ws.Range("S83").Value = 1
$ws.Range("AK83").Value = 1
$ws.Range("B82:AW82").ClearContents()
$ws.Range("A82").Value = 'Portage Park'
$ws.Range("E82").Value = 3
$ws.Range("G82").Value = 2
$ws.Range("K82").Value = 1
$ws.Range("M82").Value = 3
$ws.Range("O82").Value = 1
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = 1
$ws.Range("S82").Value = 1
$ws.Range("W82").Value = 3
$ws.Range("AC82").Value = 1
$ws.Range("AD82").Value = 1
$ws.Range("AE82").Value = 3
$ws.Range("AG82").Value = 1
$ws.Range("AI82").Value = 2
$ws.Range("AN82").Value = 2
$ws.Range("B81:AW81").ClearContents()
$ws.Range("A81").Value = 'Old Town'
$ws.Range("G81").Value = 2
$ws.Range("M81").Value = 2
$ws.Range("R81").Value = 2
$ws.Range("X81").Value = 1
$ws.Range("AG81").Value = 1
$ws.Range("B80:AW80").ClearContents()
$ws.Range("A80").Value = 'Oakland'
$ws.Range("B80").Value = 1
$ws.Range("E80").Value = 1
$ws.Range("M80").Value = 1
$ws.Range("O80").Value = 1
$ws.Range("U80").Value = 1
$ws.Range("AB80").Value = 1
$ws.Range("AC80").Value = 1
$ws.Range("AG80").Value = 1
$ws.Range("B79:AW79").ClearContents()
$ws.Range("A79").Value = 'O''Hare'
$ws.Range("Q79").Value = 1
$ws.Range("B78:AW78").ClearContents()
$ws.Range("A78").Value = 'Norwood Park'
$ws.Range("E78").Value = 1
$ws.Range("L78").Value = 1
$ws.Range("B77:AW77").ClearContents()
$ws.Range("A77").Value = 'North Park'
$ws.Range("F77").Value = 1
$ws.Range("G77").Value = 2
$ws.Range("I77").Value = 1
$ws.Range("R77").Value = 1
$ws.Range("Y77").Value = 1
$ws.Range("AA77").Value = 2
$ws.Range("AB77").Value = 1
$ws.Range("AG77").Value = 1
$ws.Range("B76:AW76").ClearContents()
$ws.Range("A76").Value = 'North Center'
$ws.Range("E76").Value = 2
$ws.Range("F76").Value = 1
$ws.Range("N76").Value = 1
$ws.Range("AG76").Value = 1
$ws.Range("AQ76").Value = 1
$ws.Range("B75:AW75").ClearContents()
$ws.Range("A75").Value = 'New City'
$ws.Range("D75").Value = 1
$ws.Range("E75").Value = 5
$ws.Range("F75").Value = 2
$ws.Range("G75").Value = 10
$ws.Range("H75").Value = 1
$ws.Range("I75").Value = 4
$ws.Range("J75").Value = 1
$ws.Range("L75").Value = 2
$ws.Range("M75").Value = 3
$ws.Range("N75").Value = 1
$ws.Range("O75").Value = 1
$ws.Range("P75").Value = 1
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = 1
$ws.Range("S75").Value = 1
$ws.Range("Y75").Value = 1
$ws.Range("AA75").Value = 2
$ws.Range("AB75").Value = 1
$ws.Range("AF75").Value = 1
$ws.Range("AI75").Value = 1
$ws.Range("AJ75").Value = 1
$ws.Range("AK75").Value = 1
$ws.Range("AM75").Value = 2
$ws.Range("AO75").Value = 3
$ws.Range("AP75").Value = 2
$ws.Range("AS75").Value = 1
$ws.Range("AT75").Value = 1
$ws.Range("AU75").Value = 1
$ws.Range("B74:AW74").ClearContents()
$ws.Range("A74").Value = 'Museum Campus'
$ws.Range("M74").Value = 1
$ws.Range("B73:AW73").ClearContents()
$ws.Range("A73").Value = 'Montclare'
$ws.Range("D73").Value = 2
$ws.Range("J73").Value = 2
$ws.Range("B72:AW72").ClearContents()
$ws.Range("A72").Value = 'Millenium Park'
$ws.Range("D72").Value = 1
$ws.Range("J72").Value = 1
$ws.Range("B71:AW71").ClearContents()
$ws.Range("A71").Value = 'Lower West Side'
$ws.Range("E71").Value = 1
$ws.Range("F71").Value = 2
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 1
$ws.Range("J71").Value = 1
$ws.Range("K71").Value = 1
$ws.Range("M71").Value = 2
$ws.Range("O71").Value = 1
$ws.Range("Y71").Value = 1
$ws.Range("Z71").Value = 1
$ws.Range("AA71").Value = 1
$ws.Range("AJ71").Value = 2
$ws.Range("AO71").Value = 1
$ws.Range("B70:AW70").ClearContents()
$ws.Range("A70").Value = 'Loop'
$ws.Range("B70").Value = 1
$ws.Range("D70").Value = 4
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 1
$ws.Range("G70").Value = 1
$ws.Range("I70").Value = 2
$ws.Range("J70").Value = 1
$ws.Range("L70").Value = 1
$ws.Range("M70").Value = 1
$ws.Range("O70").Value = 1
$ws.Range("R70").Value = 1
$ws.Range("S70").Value = 1
$ws.Range("AD70").Value = 1
$ws.Range("AE70").Value = 1
$ws.Range("AF70").Value = 1
$ws.Range("AI70").Value = 1
$ws.Range("B69:AW69").ClearContents()
$ws.Range("A69").Value = 'Lincoln Square'
$ws.Range("J69").Value = 1
$ws.Range("K69").Value = 1
$ws.Range("M69").Value = 1
$ws.Range("O69").Value = 1
$ws.Range("AB69").Value = 1
$ws.Range("AM69").Value = 1
$ws.Range("AU69").Value = 1
$ws.Range("B68:AW68").ClearContents()
$ws.Range("A68").Value = 'Kenwood'
$ws.Range("D68").Value = 3
$ws.Range("E68").Value = 1
$ws.Range("F68").Value = 5
$ws.Range("G68").Value = 3
$ws.Range("I68").Value = 3
$ws.Range("J68").Value = 2
$ws.Range("K68").Value = 1
$ws.Range("L68").Value = 3
$ws.Range("M68").Value = 8
$ws.Range("S68").Value = 1
$ws.Range("U68").Value = 1
$ws.Range("V68").Value = 1
$ws.Range("Z68").Value = 1
$ws.Range("AC68").Value = 2
$ws.Range("AD68").Value = 1
$ws.Range("AE68").Value = 2
$ws.Range("AI68").Value = 1
$ws.Range("AT68").Value = 1
$ws.Range("B67:AW67").ClearContents()
$ws.Range("A67").Value = 'Jefferson Park'
$ws.Range("D67").Value = 1
$ws.Range("F67").Value = 1
$ws.Range("M67").Value = 3
$ws.Range("AN67").Value = 1
$ws.Range("B66:AW66").ClearContents()
$ws.Range("A66").Value = 'Jackson Park'
$ws.Range("B66").Value = 1
$ws.Range("I66").Value = 1
$ws.Range("O66").Value = 1
$ws.Range("AA66").Value = 1
$ws.Range("AD66").Value = 1
$ws.Range("AG66").Value = 1
$ws.Range("AJ66").Value = 1
$ws.Range("AK66").Value = 1
$ws.Range("B65:AW65").ClearContents()
$ws.Range("A65").Value = 'Irving Park'
$ws.Range("E65").Value = 1
$ws.Range("G65").Value = 1
$ws.Range("L65").Value = 1
$ws.Range("M65").Value = 1
$ws.Range("N65").Value = 2
$ws.Range("P65").Value = 1
$ws.Range("W65").Value = 1
$ws.Range("X65").Value = 1
$ws.Range("Z65").Value = 1
$ws.Range("AE65").Value = 3
$ws.Range("AI65").Value = 1
$ws.Range("AK65").Value = 3
$ws.Range("AP65").Value = 1
$ws.Range("AW65").Value = 1
$ws.Range("B64:AW64").ClearContents()
$ws.Range("A64").Value = 'Hermosa'
$ws.Range("F64").Value = 1
$ws.Range("G64").Value = 1
$ws.Range("K64").Value = 1
$ws.Range("O64").Value = 1
$ws.Range("P64").Value = 1
$ws.Range("Q64").Value = 1
$ws.Range("U64").Value = 1
$ws.Range("X64").Value = 1
$ws.Range("AB64").Value = 1
$ws.Range("AH64").Value = 1
$ws.Range("AI64").Value = 1
$ws.Range("AJ64").Value = 1
$ws.Range("AP64").Value = 3
$ws.Range("B63:AW63").ClearContents()
$ws.Range("A63").Value = 'Hegewisch'
$ws.Range("I63").Value = 1
$ws.Range("P63").Value = 1
$ws.Range("AK63").Value = 1
$ws.Range("B62:AW62").ClearContents()
$ws.Range("A62").Value = 'Greektown'
$ws.Range("O62").Value = 1
$ws.Range("AB62").Value = 1

# New row 61: Grant Park (copy column-A style from row 60 first)
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("B61:AW61").ClearContents()
$ws.Range("A61").Value = "Grant Park"
$ws.Range("H61").Value = 1
